# Adds a new "2022-Q4" quarterly sheet (inserted right after "总计"),
# populates it with fund-holding data, and inserts a corresponding summary
# row at the top of the "总计" table (shifting the existing rows down).

function Set-TextCell($cell, $val) {
    # Forces a numeric-looking string (e.g. a fund code like "005434" or a
    # decimal like "6.02") to be stored as TEXT rather than being silently
    # coerced to a number by Excel's type inference, while leaving the
    # cell's style back at the sheet default (no explicit style index).
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $total)
$q4.Name = "2022-Q4"

# Header row
$q4.Cells.Item(1, 2).Value2 = "基金代码"
$q4.Cells.Item(1, 3).Value2 = "基金名称"
$q4.Cells.Item(1, 4).Value2 = "基金规模"
$q4.Cells.Item(1, 5).Value2 = "股票总仓位"
$q4.Cells.Item(1, 6).Value2 = "仓位占比"
$q4.Cells.Item(1, 7).Value2 = "持有市值(亿元)"
$q4.Cells.Item(1, 8).Value2 = "仓位排名"

# Reuse the header style (s="2") straight from the "总计" header cell, by
# copying formats only, so we do not depend on a Normal/Style round-trip.
$total.Cells.Item(1, 2).Copy()
for ($col = 2; $col -le 8; $col++) {
    $q4.Cells.Item(1, $col).PasteSpecial(-4122)
}

# Data rows: [code, name, scale, totalPos, posPct, heldValue, rank]
$data = @(
    @("005434", "鹏华睿投灵活配置混合A", "6.02", "82.78", "2.36", "0.1421", 5),
    @("014155", "国泰君安中证500指数增强A", "7.70", "92.93", "1.09", "0.0839", 5),
    @("014156", "国泰君安中证500指数增强C", "4.81", "92.93", "1.09", "0.0524", 5),
    @("006729", "万家中证500指数增强A", "3.13", "93.56", "1.10", "0.0344", 9),
    @("006730", "万家中证500指数增强C", "2.38", "93.56", "1.10", "0.0262", 9),
    @("015453", "中欧中证500指数增强A", "1.10", "91.84", "1.55", "0.0170", 10),
    @("005140", "华夏睿磐泰荣混合A", "3.52", "20.27", "0.24", "0.0084", 10),
    @("005141", "华夏睿磐泰荣混合C", "3.04", "20.27", "0.24", "0.0073", 10),
    @("015454", "中欧中证500指数增强C", "0.32", "91.84", "1.55", "0.0050", 10),
    @("016950", "鹏华睿投灵活配置混合C", "0.16", "82.78", "2.36", "0.0038", 5)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $q4.Cells.Item($r, 1).Value2 = $i
    Set-TextCell $q4.Cells.Item($r, 2) $row[0]
    $q4.Cells.Item($r, 3).Value2 = $row[1]
    Set-TextCell $q4.Cells.Item($r, 4) $row[2]
    Set-TextCell $q4.Cells.Item($r, 5) $row[3]
    Set-TextCell $q4.Cells.Item($r, 6) $row[4]
    Set-TextCell $q4.Cells.Item($r, 7) $row[5]
    $q4.Cells.Item($r, 8).Value2 = $row[6]
}

# Give column A on the data rows the same index-column style (s="2") used
# throughout the workbook, copied from the "总计" sheet's own A-column cell.
$total.Cells.Item(2, 1).Copy()
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $q4.Cells.Item($r, 1).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# 2) Insert a new top row in "总计" for 2022-Q4 (10 holdings, 0.38亿元),
#    shifting the existing history rows down by one.
# ---------------------------------------------------------------------
for ($r = 8; $r -ge 2; $r--) {
    $b = $total.Cells.Item($r, 2).Value2
    $c = $total.Cells.Item($r, 3).Value2
    $d = $total.Cells.Item($r, 4).Value2

    $total.Cells.Item($r, 1).Copy()
    $total.Cells.Item($r + 1, 1).PasteSpecial(-4122)

    $total.Cells.Item($r + 1, 1).Value2 = ($r - 1)
    $total.Cells.Item($r + 1, 2).Value2 = $b
    $total.Cells.Item($r + 1, 3).Value2 = $c
    $total.Cells.Item($r + 1, 4).Value2 = $d
}

$total.Cells.Item(2, 1).Value2 = 0
$total.Cells.Item(2, 2).Value2 = "2022-Q4"
$total.Cells.Item(2, 3).Value2 = 10
$total.Cells.Item(2, 4).Value2 = 0.38

# Restore the originally-active sheet (the last sheet, "2020-Q4"), since
# adding the new worksheet above shifted which tab is active.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
